$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reword the 15.1.2 indicator description stored in B4.
$ws.Range("B4").Value = "15.1.2 Доля важных с точки зрения биологического разнообразия районов суши и пресноводных районов, находящихся под охраной, в разбивке по видам экосистем"

# Move the active selection to B4, matching the state captured when the
# file was re-saved.
$ws.Range("B4").Select()
